$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 86, shifting existing rows 86-128 down to 87-129
$ws.Rows.Item(86).Insert()

# Fill in the new row 86 with the new record
$ws.Cells.Item(86, 1).Value = 5
$ws.Cells.Item(86, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(86, 3).Value = "Maule"
$ws.Cells.Item(86, 4).Value = 44489
$ws.Cells.Item(86, 5).Value = 7
$ws.Cells.Item(86, 6).Value = 100112017
$ws.Cells.Item(86, 7).Value = "Apio"
$ws.Cells.Item(86, 8).Value = "Americana (o)"
$ws.Cells.Item(86, 9).Value = "Primera"
$ws.Cells.Item(86, 10).Value = 500
$ws.Cells.Item(86, 11).Value = 7000
$ws.Cells.Item(86, 12).Value = 7000
$ws.Cells.Item(86, 13).Value = 7000
$ws.Cells.Item(86, 14).Value = "$/docena de matas"
$ws.Cells.Item(86, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(86, 16).Value = 1167
$ws.Cells.Item(86, 17).Value = 6
$ws.Cells.Item(86, 18).Value = "Hortaliza"
